$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The corrector (Furkan Kara) has now completed his pass over the review
# sheet: his name replaces the "[Developer Name]" placeholder in the
# Corrector header cell, and the Corrector column (C) is filled in with a
# result for every checklist item, mirroring the Reviewer column (B) in
# most rows.
$ws.Range("C2").Value = "Furkan Kara"

$ws.Range("B3").Copy($ws.Range("C3"))
$ws.Range("B4").Copy($ws.Range("C4"))
$ws.Range("B6").Copy($ws.Range("C6"))
$ws.Range("B7").Copy($ws.Range("C7"))
$ws.Range("B8").Copy($ws.Range("C8"))
$ws.Range("B9").Copy($ws.Range("C9"))
$ws.Range("B10").Copy($ws.Range("C10"))
$ws.Range("B12").Copy($ws.Range("C12"))
$ws.Range("B13").Copy($ws.Range("C13"))
$ws.Range("B14").Copy($ws.Range("C14"))
$ws.Range("B15").Copy($ws.Range("C15"))
$ws.Range("B16").Copy($ws.Range("C16"))

# Two rows where the corrector's finding differs from the reviewer's
# original "Rejected" verdict.
$ws.Range("F2").Copy($ws.Range("C5"))
$ws.Range("F3").Copy($ws.Range("C11"))

# Refresh the view: scroll back to the top, zoom out a bit, and leave the
# selection where Furkan finished working.
$excel.ActiveWindow.Zoom = 57
$ws.Range("C20").Select() | Out-Null
